$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 76.90000000000001
$ws.Range("I9").Value = 79.833336
$ws.Range("K9").Value = 79.833336
$ws.Range("M9").Value = 89.166664

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1955.6
$ws.Range("I40").Value = 1127.1428
$ws.Range("J40").Value = 2277.7778
$ws.Range("K40").Value = 1127.1428
$ws.Range("L40").Value = 2277.7778
$ws.Range("M40").Value = -952.1428000000001
$ws.Range("N40").Value = -2627.7778

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 3600
$ws.Range("I55").Value = 3600
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 3600
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -3386
$ws.Range("N55").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3292.3076
$ws.Range("I64").Value = 3500
$ws.Range("J64").Value = 2960
$ws.Range("K64").Value = 3500
$ws.Range("L64").Value = 2960
$ws.Range("M64").Value = -3252
$ws.Range("N64").Value = -3456

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3292.3076
$ws.Range("I67").Value = 3500
$ws.Range("J67").Value = 2960
$ws.Range("K67").Value = 3500
$ws.Range("L67").Value = 2960
$ws.Range("M67").Value = -2642
$ws.Range("N67").Value = -4676

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3685.5715
$ws.Range("I116").Value = 3808.5833
$ws.Range("K116").Value = 3808.5833
$ws.Range("M116").Value = -366.5832999999998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2409.7917
$ws.Range("I138").Value = 1021.19446
$ws.Range("J138").Value = 3242.95
$ws.Range("K138").Value = 3063.58338
$ws.Range("L138").Value = 9728.849999999999
$ws.Range("M138").Value = 2076.41662
$ws.Range("N138").Value = -20008.85

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 25500
$ws.Range("I41").Value = 1000
$ws.Range("K41").Value = 1000
$ws.Range("M41").Value = -586

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2098.524
$ws.Range("I63").Value = 2098.524
$ws.Range("K63").Value = 2098.524
$ws.Range("M63").Value = -1412.524

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2098.524
$ws.Range("I66").Value = 2098.524
$ws.Range("K66").Value = 10492.62
$ws.Range("M66").Value = -7060.619999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2805.3845
$ws.Range("I102").Value = 2788.25
$ws.Range("K102").Value = 2788.25
$ws.Range("M102").Value = -1166.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 18750
$ws.Range("I105").Value = 15000
$ws.Range("J105").Value = 30000
$ws.Range("K105").Value = 15000
$ws.Range("L105").Value = 30000
$ws.Range("M105").Value = -13253
$ws.Range("N105").Value = -33494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1101.5714
$ws.Range("I22").Value = 1288.1818
$ws.Range("J22").Value = 417.33334
$ws.Range("K22").Value = 1288.1818
$ws.Range("L22").Value = 417.33334
$ws.Range("M22").Value = -938.1818000000001
$ws.Range("N22").Value = -1117.33334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 3268138
$ws.Range("I2").Value = 4902081
$ws.Range("J2").Value = 252
$ws.Range("K2").Value = 29412486
$ws.Range("L2").Value = 1512
$ws.Range("M2").Value = -29412373
$ws.Range("N2").Value = -1738

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5312.5
$ws.Range("I80").Value = 2625
$ws.Range("J80").Value = 8000
$ws.Range("K80").Value = 7875
$ws.Range("L80").Value = 24000
$ws.Range("M80").Value = -6939
$ws.Range("N80").Value = -25872

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 5312.5
$ws.Range("I83").Value = 2625
$ws.Range("J83").Value = 8000
$ws.Range("K83").Value = 23625
$ws.Range("L83").Value = 72000
$ws.Range("M83").Value = -18945
$ws.Range("N83").Value = -81360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 795
$ws.Range("I121").Value = 172
$ws.Range("J121").Value = 1833.3334
$ws.Range("K121").Value = 516
$ws.Range("L121").Value = 5500.0002
$ws.Range("M121").Value = 794
$ws.Range("N121").Value = -8120.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 4324.375
$ws.Range("I141").Value = 1749.091
$ws.Range("K141").Value = 5247.272999999999
$ws.Range("M141").Value = -67.27299999999923

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 44208.332
$ws.Range("J135").Value = 44208.332
$ws.Range("L135").Value = 44208.332
$ws.Range("N135").Value = -54348.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1299.7142
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1419.6
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1419.6
$ws.Range("N22").Value = -2009.6
$ws.Range("M22").Value = -705

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1299.7142
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 1419.6
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1419.6
$ws.Range("N27").Value = -1633.6
$ws.Range("M27").Value = -893

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3063.6365
$ws.Range("I46").Value = 3100
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 3100
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -2912
$ws.Range("N46").Value = -3376

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 398.8
$ws.Range("I55").Value = 398
$ws.Range("J55").Value = 400
$ws.Range("K55").Value = 398
$ws.Range("L55").Value = 400
$ws.Range("M55").Value = -225
$ws.Range("N55").Value = -746

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 16666.666
$ws.Range("J109").Value = 16666.666
$ws.Range("L109").Value = 16666.666
$ws.Range("N109").Value = -19440.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3971.4285
$ws.Range("I62").Value = 4900
$ws.Range("J62").Value = 3600
$ws.Range("K62").Value = 4900
$ws.Range("L62").Value = 3600
$ws.Range("M62").Value = -4276
$ws.Range("N62").Value = -4848

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 3971.4285
$ws.Range("I65").Value = 4900
$ws.Range("J65").Value = 3600
$ws.Range("K65").Value = 24500
$ws.Range("L65").Value = 18000
$ws.Range("M65").Value = -21380
$ws.Range("N65").Value = -24240
